# CareDesk appointments log — add call-transfer/reschedule functionality.
#
# Rows 2 (Hassan) and 3 (Sunny) were originally logged as fresh "book"
# actions. The bot now reschedules them to a new preferred slot and keeps
# a record of the previously booked slot in the "Existing Appointment"
# column (K). A brand-new booking (Manju) is appended as row 4.
#
# Date-/number-looking text (dates, phone numbers, ages) is written into
# cells pre-formatted as Text ("@") so Excel's automatic type detection
# does not silently convert them to real dates/numbers, matching how
# these values were originally stored (as plain text) in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Hassan -> rescheduled --------------------------------------
$ws.Range("B2").Value = "reschedule"
$ws.Range("F2").Value = "Miss Khalil, orthopedics"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2020-02-20"

$ws.Range("I2").Value = "morning"
$ws.Range("K2").Value = "2026-01-21 evening"

# --- Row 3: Sunny -> rescheduled ----------------------------------------
$ws.Range("B3").Value = "reschedule"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2026-01-25"

$ws.Range("I3").Value = "evening"
$ws.Range("K3").Value = "2026-01-21 morning"
$ws.Range("L3").ClearContents()

# --- Row 4: Manju -> new booking appended --------------------------------
$ws.Range("A4").Value = "2026-01-21 17:37:10"
$ws.Range("B4").Value = "book"
$ws.Range("C4").Value = "Manju"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "22"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "01518307641"

$ws.Range("F4").Value = "Dr. Rachel Morgan, Family Medicine"
$ws.Range("G4").Value = "cold issue"
$ws.Range("H4").Value = "Monday"
$ws.Range("I4").Value = "9 AM"
$ws.Range("J4").Value = "first visit"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
